$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> (expected old text, new text)
# Cell text is replaced directly (instead of a global Find/Replace) to
# avoid collisions, since some 'new' values equal other cells' 'old'
# values (e.g. 84÷6= -> 72÷8= while a different cell already holds 72÷8=).
$replacements = @(
    @{ Row = 1; Col = 1; Old = "85÷4="; New = "86÷9=" }
    @{ Row = 1; Col = 2; Old = "53÷7="; New = "85÷7=" }
    @{ Row = 1; Col = 3; Old = "44÷6="; New = "12÷5=" }
    @{ Row = 1; Col = 4; Old = "73÷6="; New = "81÷2=" }
    @{ Row = 1; Col = 5; Old = "89÷9="; New = "21÷5=" }
    @{ Row = 5; Col = 1; Old = "86÷7="; New = "87÷3=" }
    @{ Row = 5; Col = 2; Old = "58÷8="; New = "93÷5=" }
    @{ Row = 5; Col = 3; Old = "31÷4="; New = "26÷7=" }
    @{ Row = 5; Col = 4; Old = "57÷3="; New = "36÷3=" }
    @{ Row = 5; Col = 5; Old = "84÷6="; New = "72÷8=" }
    @{ Row = 9; Col = 1; Old = "99÷5="; New = "77÷4=" }
    @{ Row = 9; Col = 2; Old = "93÷3="; New = "16÷8=" }
    @{ Row = 9; Col = 3; Old = "93÷9="; New = "19÷8=" }
    @{ Row = 9; Col = 4; Old = "98÷2="; New = "18÷4=" }
    @{ Row = 9; Col = 5; Old = "40÷8="; New = "69÷5=" }
    @{ Row = 13; Col = 1; Old = "34÷8="; New = "67÷2=" }
    @{ Row = 13; Col = 2; Old = "89÷3="; New = "54÷5=" }
    @{ Row = 13; Col = 3; Old = "54÷7="; New = "75÷2=" }
    @{ Row = 13; Col = 4; Old = "37÷8="; New = "94÷2=" }
    @{ Row = 13; Col = 5; Old = "33÷7="; New = "60÷5=" }
    @{ Row = 17; Col = 1; Old = "72÷8="; New = "41÷9=" }
    @{ Row = 17; Col = 2; Old = "14÷6="; New = "30÷2=" }
    @{ Row = 17; Col = 3; Old = "72÷3="; New = "67÷2=" }
    @{ Row = 17; Col = 4; Old = "43÷6="; New = "50÷2=" }
    @{ Row = 17; Col = 5; Old = "51÷9="; New = "36÷5=" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $r = $cell.Range
    # Trim the trailing end-of-cell marker so we only touch the visible text
    $r.End = $r.End - 1
    if ($r.Text -ne $item.Old) {
        Write-Host "WARNING: cell ($($item.Row),$($item.Col)) expected $($item.Old) but found $($r.Text)"
    }
    $r.Text = $item.New
}

Write-Host "Done applying replacements."
